$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal text, avoiding Excel
# auto-converting number-looking strings (e.g. "563.70" -> 563.7).
function Set-TextValue {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "63.685.94"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "3.145.34"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws "D5" "563.70"
$ws.Range("E5").Value = "  -3.15%  "
Set-TextValue $ws "D6" "166.44"
$ws.Range("E6").Value = "  -6.90%  "
Set-TextValue $ws "D7" "0.600"
$ws.Range("E7").Value = "  -6.47%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.150.18"
$ws.Range("E9").Value = "  -3.69%  "
Set-TextValue $ws "D10" "0.119"
$ws.Range("E10").Value = "  -4.17%  "
Set-TextValue $ws "D11" "6.64"
$ws.Range("E11").Value = "  -1.10%  "
Set-TextValue $ws "D12" "0.380"
$ws.Range("E12").Value = "  -5.05%  "
$ws.Range("D13").Value = "3.699.55"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D15").Value = "63.790.83"
$ws.Range("E15").Value = "  -2.96%  "
Set-TextValue $ws "D16" "25.06"
$ws.Range("E16").Value = "  -3.28%  "
Set-TextValue $ws "D17" "0.0000156"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "3.156.90"
$ws.Range("E18").Value = "  -3.11%  "
Set-TextValue $ws "D19" "409.55"
$ws.Range("E19").Value = "  -3.95%  "
Set-TextValue $ws "D20" "12.67"
$ws.Range("E20").Value = "  -3.75%  "
Set-TextValue $ws "D21" "5.28"
$ws.Range("E21").Value = "  -3.33%  "
Set-TextValue $ws "D22" "7.04"
$ws.Range("E22").Value = "  -4.07%  "
Set-TextValue $ws "D23" "0.997"
$ws.Range("E23").Value = "  -0.29%  "
Set-TextValue $ws "D24" "70.22"
$ws.Range("E24").Value = "  -2.25%  "
Set-TextValue $ws "D25" "0.199"
$ws.Range("E25").Value = "  +1.91%  "
Set-TextValue $ws "D26" "0.486"
$ws.Range("E26").Value = "  -4.10%  "
Set-TextValue $ws "D27" "0.0000105"
$ws.Range("E27").Value = "  -6.35%  "
Set-TextValue $ws "D28" "8.61"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("E29").Value = "  +0.07%  "
Set-TextValue $ws "D30" "1.82"
$ws.Range("E30").Value = "  -6.83%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D31" "0.999"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D32" "21.53"
$ws.Range("E32").Value = "  -2.77%  "
Set-TextValue $ws "D33" "4.93"
$ws.Range("E33").Value = "  -3.58%  "
Set-TextValue $ws "D34" "6.27"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  -5.41%  "
Set-TextValue $ws "D36" "155.09"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").Value = "2.704.20"
$ws.Range("E38").Value = "  -2.73%  "
Set-TextValue $ws "D39" "1.68"
$ws.Range("E39").Value = "  -6.36%  "
Set-TextValue $ws "D40" "24.39"
$ws.Range("E40").Value = "  -7.04%  "
Set-TextValue $ws "D41" "4.12"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D42" "38.60"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D43" "0.705"
$ws.Range("E43").Value = "  -7.79%  "
Set-TextValue $ws "D44" "0.0617"
$ws.Range("E44").Value = "  -5.81%  "
Set-TextValue $ws "D45" "5.55"
$ws.Range("E45").Value = "  -5.85%  "
Set-TextValue $ws "D46" "0.0259"
$ws.Range("E46").Value = "  -2.63%  "
Set-TextValue $ws "D47" "21.55"
$ws.Range("E47").Value = "  -6.04%  "
Set-TextValue $ws "D48" "291.65"
$ws.Range("E48").Value = "  -7.07%  "
Set-TextValue $ws "D49" "2.01"
$ws.Range("E49").Value = "  -11.47%  "
Set-TextValue $ws "D50" "1.00"
$ws.Range("E50").Value = "  -0.02%  "
Set-TextValue $ws "D51" "0.0985"
$ws.Range("E51").Value = "  -5.46%  "
